$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price and Volume columns for data rows (2-51) are treated as text,
# matching the source workbook where these values are stored as inline strings
# (not numbers), e.g. "1.001", "22.417.57", "  -0.11%  ".
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '22.417.57'
$ws.Range("E2").Value = '  -0.11%  '
$ws.Range("D3").Value = '1.564.35'
$ws.Range("E3").Value = '  -0.50%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.23%  '
$ws.Range("D5").Value = '1.001'
$ws.Range("E5").Value = '  -0.23%  '
$ws.Range("D6").Value = '287.16'
$ws.Range("E6").Value = '  -1.30%  '
$ws.Range("D7").Value = '0.3633'
$ws.Range("E7").Value = '  -3.04%  '
$ws.Range("D8").Value = '49.78'
$ws.Range("E8").Value = '  +0.18%  '
$ws.Range("D9").Value = '0.3351'
$ws.Range("E9").Value = '  -1.52%  '
$ws.Range("D10").Value = '1.126'
$ws.Range("E10").Value = '  -1.97%  '
$ws.Range("D11").Value = '0.07406'
$ws.Range("E11").Value = '  -1.87%  '
$ws.Range("D12").Value = '1.001'
$ws.Range("E12").Value = '  -0.27%  '
$ws.Range("D13").Value = '20.90'
$ws.Range("E13").Value = '  -1.85%  '
$ws.Range("D14").Value = '5.928'
$ws.Range("E14").Value = '  -1.23%  '
$ws.Range("D15").Value = '6.880'
$ws.Range("E15").Value = '  -1.00%  '
$ws.Range("D16").Value = '1.564.40'
$ws.Range("E16").Value = '  -1.29%  '
$ws.Range("D17").Value = '0.00001098'
$ws.Range("E17").Value = '  -1.70%  '
$ws.Range("D18").Value = '89.18'
$ws.Range("E18").Value = '  -1.84%  '
$ws.Range("D19").Value = '0.06730'
$ws.Range("E19").Value = '  -0.15%  '
$ws.Range("D20").Value = '1.001'
$ws.Range("E20").Value = '  -0.27%  '
$ws.Range("D21").Value = '6.308'
$ws.Range("E21").Value = '  +1.14%  '
$ws.Range("D22").Value = '15.99'
$ws.Range("E22").Value = '  -2.23%  '
$ws.Range("D23").Value = '11.95'
$ws.Range("E23").Value = '  -1.79%  '
$ws.Range("D24").Value = '22.406.06'
$ws.Range("E24").Value = '  -0.12%  '
$ws.Range("D25").Value = '2.380'
$ws.Range("E25").Value = '  +0.75%  '
$ws.Range("D26").Value = '2.535'
$ws.Range("E26").Value = '  -2.97%  '
$ws.Range("D27").Value = '149.14'
$ws.Range("E27").Value = '  +0.67%  '
$ws.Range("D28").Value = '19.58'
$ws.Range("E28").Value = '  -2.73%  '
$ws.Range("D29").Value = '4.995'
$ws.Range("E29").Value = '  -0.03%  '
$ws.Range("D30").Value = '123.04'
$ws.Range("E30").Value = '  -2.21%  '
$ws.Range("D31").Value = '1.739.74'
$ws.Range("E31").Value = '  -1.01%  '
$ws.Range("D32").Value = '1.057'
$ws.Range("E32").Value = '  +3.17%  '
$ws.Range("D33").Value = '6.096'
$ws.Range("E33").Value = '  -0.15%  '
$ws.Range("D34").Value = '1.974'
$ws.Range("E34").Value = '  -0.70%  '
$ws.Range("D35").Value = '9.519'
$ws.Range("E35").Value = '  -3.36%  '
$ws.Range("D36").Value = '0.08254'
$ws.Range("E36").Value = '  -2.12%  '
$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").Value = '0.02384'
$ws.Range("E37").Value = '  -3.01%  '
$ws.Range("B38").Value = 'TrustWalletToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D38").Value = '1.310'
$ws.Range("E38").Value = '  -4.49%  '
$ws.Range("D39").Value = '0.2211'
$ws.Range("E39").Value = '  -3.12%  '
$ws.Range("D40").Value = '0.06343'
$ws.Range("E40").Value = '  -3.31%  '
$ws.Range("D41").Value = '5.310'
$ws.Range("E41").Value = '  -2.69%  '
$ws.Range("D42").Value = '11.21'
$ws.Range("E42").Value = '  -1.47%  '
$ws.Range("B43").Value = 'TheSandbox'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D43").Value = '0.6069'
$ws.Range("E43").Value = '  -3.44%  '
$ws.Range("B44").Value = 'Frax'
$ws.Range("C44").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D44").Value = '1.001'
$ws.Range("E44").Value = '  -0.21%  '
$ws.Range("D45").Value = '13.73'
$ws.Range("E45").Value = '  -2.31%  '
$ws.Range("D46").Value = '3.767'
$ws.Range("E46").Value = '  -1.36%  '
$ws.Range("D47").Value = '0.5705'
$ws.Range("E47").Value = '  -2.60%  '
$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").Value = '2.009'
$ws.Range("E48").Value = '  -4.40%  '
$ws.Range("B49").Value = 'EOS'
$ws.Range("C49").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D49").Value = '1.232'
$ws.Range("E49").Value = '  +0.68%  '
$ws.Range("D50").Value = '124.20'
$ws.Range("E50").Value = '  -4.39%  '
$ws.Range("D51").Value = '0.07249'
$ws.Range("E51").Value = '  -1.17%  '
